$wb = $excel.ActiveWorkbook

# Fix the typo in the second worksheet's name:
# "Puerto Prinsesa" -> "Puerto Princesa"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Puerto Princesa"

# Make the "Puerto Princesa" sheet the active tab (previously "Calapan" was active).
$ws2.Activate()

# Update the current selection on the active sheet to a single cell E23.
$ws2.Range("E23").Select()
